$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("readings")

# New rows of readings data to append (id, consumerId, date, previousReading, currentReading, consumption, total)
$rows = @(
    @(8,  3,  "2025-06-17T10:48:25.637Z", 0,   23,  23,  11.5),
    @(9,  4,  "2025-06-17T10:48:32.842Z", 0,   52,  52,  62.4),
    @(10, 25, "2025-06-17T10:48:38.621Z", 0,   74,  74,  37),
    @(11, 21, "2025-06-17T10:48:56.958Z", 199, 207, 8,   4),
    @(12, 22, "2025-06-17T10:49:02.956Z", 0,   63,  63,  31.5),
    @(13, 26, "2025-06-17T10:49:09.165Z", 0,   56,  56,  67.2),
    @(14, 26, "2025-06-17T10:49:44.465Z", 56,  74,  18,  21.599999999999998)
)

$startRow = 9
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
}
